$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Country"
$ws.Range("B1").Value = "Number of People Living with HIV"
$ws.Range("C1").Value = "Descrimination Percent"

$ws.Range("C2").Select()
